# Update Name of Algo
# Apply updated KNN-imputed values to column C for the specified rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = -12.223
    6   = -12.318
    7   = -12.909
    16  = -13.164
    20  = -12.205
    28  = -12.848
    29  = -11.736
    32  = -13.222
    40  = -12.237
    46  = -13.827
    51  = -11.187
    52  = -11.614
    57  = -13.632
    59  = -12.918
    62  = -13.555
    66  = -10.836
    73  = -12.401
    74  = -11.982
    92  = -11.104
    100 = -12.747
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
